$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be stored as text even when the string looks like a
    # number (e.g. "226.33"), matching the original workbook's inlineStr cells.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "34.376.00"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.787.80"
$ws.Range("E3").Value = "  +0.33%  "

# Row 5 - BNB
Set-TextValue "D5" "226.33"
$ws.Range("E5").Value = "  +0.23%  "

# Row 6 - XRP
Set-TextValue "D6" "0.555"
$ws.Range("E6").Value = "  +1.62%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - Solana
Set-TextValue "D8" "32.58"
$ws.Range("E8").Value = "  +1.51%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.68%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0689"
$ws.Range("E10").Value = "  +0.48%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0946"

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("E12").Value = "  +0.42%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.791.59"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14 - Chainlink
Set-TextValue "D14" "11.03"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.70%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "34.379.62"
$ws.Range("E16").Value = "  +0.79%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +2.40%  "

# Row 18 - Litecoin
Set-TextValue "D18" "68.30"
$ws.Range("E18").Value = "  +1.08%  "

# Row 19 - swap BitcoinCash -> ShibaInu
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D19" "0.0₃0794"
$ws.Range("E19").Value = "  +0.98%  "

# Row 20 - swap ShibaInu -> BitcoinCash
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "244.58"
$ws.Range("E20").Value = "  -0.30%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  +2.90%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.07%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +0.96%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.43%  "

# Row 25 - Monero
Set-TextValue "D25" "165.74"
$ws.Range("E25").Value = "  +2.17%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +2.35%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +1.42%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.43%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("E30").Value = "  +7.23%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +1.56%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.85%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +0.22%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.91%  "

# Row 35 - RenderToken
$ws.Range("E35").Value = "  +5.35%  "

# Row 36 - Maker
Set-TextValue "D36" "1.410.46"
$ws.Range("E36").Value = "  -2.56%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +4.66%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value = "  +2.78%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.14%  "

# Row 40 - Aave
Set-TextValue "D40" "84.07"
$ws.Range("E40").Value = "  +3.28%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +0.78%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  +2.60%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  +2.86%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "13.87"
$ws.Range("E44").Value = "  +2.02%  "

# Row 45 - Kaspa
Set-TextValue "D45" "0.0525"
$ws.Range("E45").Value = "  +0.84%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  +2.85%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  +0.06%  "

# Row 48 - RocketPoolETH
Set-TextValue "D48" "1.946.39"
$ws.Range("E48").Value = "  +0.41%  "

# Row 49 - Quant
Set-TextValue "D49" "105.23"
$ws.Range("E49").Value = "  +0.56%  "

# Row 50 - PaxDollar
$ws.Range("E50").Value = "  -0.14%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  -2.56%  "
